$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.221.33"
$ws.Range("E2").Value = "  -1.38%  "

$ws.Range("D3").Value = "2.187.82"
$ws.Range("E3").Value = "  -5.84%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "'295.00"
$ws.Range("E5").Value = "  -3.61%  "

$ws.Range("D6").Value = "'81.01"
$ws.Range("E6").Value = "  -3.25%  "

$ws.Range("E7").Value = "  -3.24%  "

$ws.Range("E8").Value = "  +0.17%  "

$ws.Range("D9").Value = "'0.464"
$ws.Range("E9").Value = "  -3.73%  "

$ws.Range("E10").Value = "  -5.96%  "

$ws.Range("D11").Value = "'28.92"
$ws.Range("E11").Value = "  -3.02%  "

$ws.Range("D12").Value = "'47.03"
$ws.Range("E12").Value = "  -10.18%  "

$ws.Range("E13").Value = "  -2.16%  "

$ws.Range("D14").Value = "2.529.82"
$ws.Range("E14").Value = "  -5.58%  "

$ws.Range("D15").Value = "'6.21"
$ws.Range("E15").Value = "  -2.42%  "

$ws.Range("D16").Value = "'13.88"
$ws.Range("E16").Value = "  -5.54%  "

$ws.Range("D17").Value = "2.181.90"
$ws.Range("E17").Value = "  -6.06%  "

$ws.Range("D18").Value = "'0.705"
$ws.Range("E18").Value = "  -5.09%  "

$ws.Range("D19").Value = "39.112.85"
$ws.Range("E19").Value = "  -1.46%  "

$ws.Range("D20").Value = "0.0₃0866"
$ws.Range("E20").Value = "  -3.66%  "

$ws.Range("E21").Value = "  -5.78%  "

$ws.Range("D22").Value = "'64.73"
$ws.Range("E22").Value = "  -4.05%  "

$ws.Range("D23").Value = "'10.24"
$ws.Range("E23").Value = "  -2.87%  "

$ws.Range("D24").Value = "'224.99"
$ws.Range("E24").Value = "  -3.76%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("E26").Value = "  -5.33%  "

$ws.Range("E27").Value = "  +0.15%  "

$ws.Range("D28").Value = "'22.42"
$ws.Range("E28").Value = "  -3.48%  "

$ws.Range("E29").Value = "  -1.94%  "

$ws.Range("E30").Value = "  -1.22%  "

$ws.Range("D31").Value = "'149.87"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("D32").Value = "'31.65"
$ws.Range("E32").Value = "  -7.50%  "

$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("D34").Value = "'4.78"
$ws.Range("E34").Value = "  -5.92%  "

$ws.Range("E35").Value = "  -3.68%  "

$ws.Range("E36").Value = "  -4.06%  "

$ws.Range("D37").Value = "'0.109"
$ws.Range("E37").Value = "  -2.83%  "

$ws.Range("D38").Value = "'15.24"
$ws.Range("E38").Value = "  -2.55%  "

$ws.Range("D39").Value = "'0.0958"
$ws.Range("E39").Value = "  -2.75%  "

$ws.Range("E40").Value = "  -4.79%  "

$ws.Range("E41").Value = "  -3.96%  "

$ws.Range("D42").Value = "'3.57"
$ws.Range("E42").Value = "  -5.27%  "

$ws.Range("D43").Value = "1.890.40"
$ws.Range("E43").Value = "  -2.22%  "

$ws.Range("E44").Value = "  -9.34%  "

$ws.Range("D45").Value = "'0.0258"
$ws.Range("E45").Value = "  -1.88%  "

$ws.Range("D46").Value = "'16.01"
$ws.Range("E46").Value = "  -8.07%  "

$ws.Range("E47").Value = "  -2.90%  "

$ws.Range("D48").Value = "'2.58"
$ws.Range("E48").Value = "  -1.60%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.398.33"
$ws.Range("E49").Value = "  -5.82%  "

$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "'70.86"
$ws.Range("E50").Value = "  +0.73%  "

$ws.Range("D51").Value = "'86.45"
$ws.Range("E51").Value = "  -5.70%  "
